$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Duplicate the last service-log entry (row 21, card "23") down into a new
# row 22 -- this is the new event being logged, starting out as a copy of
# the same card / date / correction / servicer values, with the
# measurement columns still blank (values not recorded yet).
$ws.Range("A21:P21").Copy($ws.Range("A22:P22"))
$excel.CutCopyMode = 0

# Row 21's measurement/event columns were left blank; backfill them with
# "nan", the workbook's established placeholder for a missing text value
# (matches the fully-populated sibling entry in row 20).
$ws.Range("B21:K21").Value = "nan"
$ws.Range("M21").Value = "nan"
$ws.Range("P21").Value = "nan"
